$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that currently sits after "FTP - 003"
#    (a new one will be (re)created further down, right where the live edit
#    cursor ends up after the "Vinduet ... usecase 02." rewrite).
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2. "Hovedscenariet fortsætter fra pkt. 2." -> split/rewrite into:
#    "Vinduet fortsætter fr" + _GoBack bookmark + "a " + "usecase"(spell
#    checked) + " 02."
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Hovedscenariet fortsætter fra pkt. 2.")
$para = $rng.Paragraphs(1)
$insertPoint = $d.Range($para.Range.Start, $para.Range.Start)

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>Vinduet forts&#230;tter fr</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>usecase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> 02.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml1)

$rng = $d.Content
$rng.Find.Execute("Hovedscenariet fortsætter fra pkt. 2.")
$rng.Delete()

# ---------------------------------------------------------------------------
# 3. "Systemet kommer med [fejlmeddelse] der " -> re-typed without the
#    spell-check markers, run boundaries shifted mid-word:
#    "Systemet kommer med fejlmedde" + "le" + "lse der "
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Systemet kommer med ")
$para2 = $rng.Paragraphs(1)
$insertPoint2 = $d.Range($para2.Range.Start, $para2.Range.Start)

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>Systemet kommer med fejlmedde</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>le</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">lse der </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint2.InsertXML($xml2)

$rng = $d.Content
$rng.Find.Execute("Systemet kommer med fejlmeddelse der ")
$rng.Delete()

# ---------------------------------------------------------------------------
# 4. The list item "Systemet gemmer oplysningerne." moves from list level 1
#    to list level 2 (w:ilvl 0 -> 1) within the same numId.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Systemet gemmer oplysningerne.")
$para3 = $rng.Paragraphs(1)
$para3.Range.ListFormat.ListLevelNumber = 2

Write-Output "done"
